$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "pt_max" column (column H). Select the whole column first
# (as a user would by clicking the column header) and then delete it so
# the remaining columns I:R shift left to H:Q. Excel takes care of
# updating all formula references, shared-string usage, the sheet
# dimension and the column-width spec automatically.
$ws.Columns("H:H").Select()
$ws.Columns("H:H").Delete()
